# Refresh the cryptos price/volume snapshot (GitHub Actions daily update).
#
# Column D ("Price") and column E ("Volume(1h)") hold plain text in this sheet
# (prices use "." as a thousands separator, e.g. "25.932.01", and some entries
# use subscript-digit notation, e.g. "0.0₅xxxx"). A handful of the refreshed
# D-column values now read as plain decimals (e.g. "0.5061"); those are written
# with a leading apostrophe - exactly what typing them into Excel does - so the
# Range.Value setter keeps storing literal text instead of converting the cell to
# a number. Rows 42/43 (FraxShare / BabyDogeCoin) also swap ranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.932.01'
$ws.Range('D3').Value = '1.641.22'
$ws.Range('E3').Value = '  -0.08%  '
$ws.Range('D4').Value = "'" + '1.002'
$ws.Range('E4').Value = '  -0.32%  '
$ws.Range('D5').Value = "'" + '215.22'
$ws.Range('E5').Value = '  -0.31%  '
$ws.Range('D6').Value = "'" + '0.5061'
$ws.Range('E6').Value = '  +0.32%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = "'" + '0.2563'
$ws.Range('E8').Value = '  -0.62%  '
$ws.Range('D9').Value = "'" + '0.06371'
$ws.Range('E9').Value = '  -0.23%  '
$ws.Range('D10').Value = "'" + '19.47'
$ws.Range('E10').Value = '  -0.38%  '
$ws.Range('D11').Value = "'" + '0.07745'
$ws.Range('E11').Value = '  -0.09%  '
$ws.Range('E12').Value = '  +0.31%  '
$ws.Range('D13').Value = '1.649.31'
$ws.Range('E13').Value = '  +0.55%  '
$ws.Range('D14').Value = "'" + '0.5443'
$ws.Range('E14').Value = '  -0.31%  '
$ws.Range('D15').Value = '0.0₅7818'
$ws.Range('E15').Value = '  -1.03%  '
$ws.Range('D16').Value = "'" + '64.26'
$ws.Range('D17').Value = '25.976.75'
$ws.Range('E17').Value = '  +0.21%  '
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('D19').Value = "'" + '197.19'
$ws.Range('E19').Value = '  -2.53%  '
$ws.Range('D20').Value = "'" + '4.435'
$ws.Range('E20').Value = '  +0.90%  '
$ws.Range('D21').Value = "'" + '9.931'
$ws.Range('E21').Value = '  +0.38%  '
$ws.Range('E22').Value = '  +1.05%  '
$ws.Range('D23').Value = "'" + '1.005'
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').Value = "'" + '1.897'
$ws.Range('E24').Value = '  +1.64%  '
$ws.Range('D25').Value = "'" + '140.97'
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').Value = "'" + '0.1168'
$ws.Range('E26').Value = '  +2.80%  '
$ws.Range('D27').Value = "'" + '6.880'
$ws.Range('E27').Value = '  +1.42%  '
$ws.Range('D28').Value = "'" + '15.68'
$ws.Range('E28').Value = '  +0.12%  '
$ws.Range('D29').Value = "'" + '1.237'
$ws.Range('E29').Value = '  -0.66%  '
$ws.Range('D30').Value = "'" + '0.04966'
$ws.Range('E30').Value = '  -0.31%  '
$ws.Range('D31').Value = "'" + '3.256'
$ws.Range('E31').Value = '  -0.61%  '
$ws.Range('D32').Value = "'" + '3.184'
$ws.Range('E32').Value = '  -0.60%  '
$ws.Range('E33').Value = '  -0.95%  '
$ws.Range('D34').Value = "'" + '2.366'
$ws.Range('E34').Value = '  -0.33%  '
$ws.Range('D35').Value = "'" + '0.8941'
$ws.Range('E35').Value = '  +0.16%  '
$ws.Range('D36').Value = "'" + '2.589'
$ws.Range('E36').Value = '  -1.74%  '
$ws.Range('D37').Value = '1.133.65'
$ws.Range('E37').Value = '  -1.64%  '
$ws.Range('D38').Value = "'" + '0.5444'
$ws.Range('E38').Value = '  -2.86%  '
$ws.Range('D39').Value = "'" + '0.01559'
$ws.Range('E39').Value = '  -0.36%  '
$ws.Range('D40').Value = "'" + '1.002'
$ws.Range('E40').Value = '  -0.34%  '
$ws.Range('D41').Value = "'" + '2.540'
$ws.Range('E41').Value = '  -1.05%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = "'" + '5.578'
$ws.Range('E42').Value = '  -1.98%  '
$ws.Range('B43').Value = 'BabyDogeCoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D43').Value = '0.0₈127'
$ws.Range('E43').Value = '  +7.46%  '
$ws.Range('D44').Value = "'" + '0.8182'
$ws.Range('E44').Value = '  +1.20%  '
$ws.Range('D45').Value = "'" + '99.60'
$ws.Range('E45').Value = '  -0.20%  '
$ws.Range('D46').Value = '1.776.64'
$ws.Range('E46').Value = '  -0.12%  '
$ws.Range('D47').Value = "'" + '0.4533'
$ws.Range('E47').Value = '  -0.30%  '
$ws.Range('D48').Value = "'" + '1.002'
$ws.Range('E49').Value = '  -0.58%  '
$ws.Range('D50').Value = "'" + '0.05069'
$ws.Range('E50').Value = '  +0.17%  '
$ws.Range('E51').Value = '  +0.39%  '
